{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"22\u00d717=374\", \"61\u00d798=5978\"],\n  [\"44\u00d784=3696\", \"86\u00d714=1204\"],\n  [\"77\u00d781=6237\", \"55\u00d729=1595\"],\n  [\"20\u00d733=660\", \"60\u00d785=5100\"],\n  [\"56\u00d726=1456\", \"87\u00d735=3045\"],\n  [\"69\u00d716=1104\", \"20\u00d761=1220\"],\n  [\"44\u00d762=2728\", \"24\u00d741=984\"],\n  [\"86\u00d788=7568\", \"90\u00d785=7650\"],\n  [\"48\u00d765=3120\", \"63\u00d793=5859\"],\n  [\"26\u00d763=1638\", \"73\u00d771=5183\"],\n  [\"44\u00d764=2816\", \"61\u00d765=3965\"],\n  [\"24\u00d739=936\", \"41\u00d793=3813\"],\n  [\"19\u00d793=1767\", \"58\u00d757=3306\"],\n  [\"40\u00d777=3080\", \"87\u00d720=1740\"],\n  [\"17\u00d761=1037\", \"39\u00d718=702\"],\n  [\"70\u00d774=5180\", \"12\u00d740=480\"],\n  [\"17\u00d794=1598\", \"85\u00d726=2210\"],\n  [\"63\u00d729=1827\", \"62\u00d726=1612\"],\n  [\"61\u00d717=1037\", \"91\u00d795=8645\"],\n  [\"96\u00d795=9120\", \"94\u00d757=5358\"],\n  [\"52\u00d746=2392\", \"16\u00d748=768\"],\n  [\"75\u00d755=4125\", \"46\u00d718=828\"],\n  [\"88\u00d731=2728\", \"24\u00d750=1200\"],\n  [\"84\u00d766=5544\", \"82\u00d712=984\"],\n  [\"50\u00d795=4750\", \"98\u00d713=1274\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"22\u00d717=374\", \"61\u00d798=5978\"),\n    @(\"44\u00d784=3696\", \"86\u00d714=1204\"),\n    @(\"77\u00d781=6237\", \"55\u00d729=1595\"),\n    @(\"20\u00d733=660\", \"60\u00d785=5100\"),\n    @(\"56\u00d726=1456\", \"87\u00d735=3045\"),\n    @(\"69\u00d716=1104\", \"20\u00d761=1220\"),\n    @(\"44\u00d762=2728\", \"24\u00d741=984\"),\n    @(\"86\u00d788=7568\", \"90\u00d785=7650\"),\n    @(\"48\u00d765=3120\", \"63\u00d793=5859\"),\n    @(\"26\u00d763=1638\", \"73\u00d771=5183\"),\n    @(\"44\u00d764=2816\", \"61\u00d765=3965\"),\n    @(\"24\u00d739=936\", \"41\u00d793=3813\"),\n    @(\"19\u00d793=1767\", \"58\u00d757=3306\"),\n    @(\"40\u00d777=3080\", \"87\u00d720=1740\"),\n    @(\"17\u00d761=1037\", \"39\u00d718=702\"),\n    @(\"70\u00d774=5180\", \"12\u00d740=480\"),\n    @(\"17\u00d794=1598\", \"85\u00d726=2210\"),\n    @(\"63\u00d729=1827\", \"62\u00d726=1612\"),\n    @(\"61\u00d717=1037\", \"91\u00d795=8645\"),\n    @(\"96\u00d795=9120\", \"94\u00d757=5358\"),\n    @(\"52\u00d746=2392\", \"16\u00d748=768\"),\n    @(\"75\u00d755=4125\", \"46\u00d718=828\"),\n    @(\"88\u00d731=2728\", \"24\u00d750=1200\"),\n    @(\"84\u00d766=5544\", \"82\u00d712=984\"),\n    @(\"50\u00d795=4750\", \"98\u00d713=1274\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$old, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, 1, [ref]$false, $new, 2)\n}\n"}
